$d = $word.ActiveDocument
$d.Content.Find.Execute("3 October 2018", $true, $false, $false, $false, $false, $true, 1, $false, "6 October 2018", 2)
